# Update the "updated at" timestamp stamped on every data row of the
# two stats sheets ("Главные" and "Линейные"). The value lives in
# column AA (rows 2-26) and must stay a plain text string, not get
# reinterpreted as a date by Excel.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-29 07:04:49"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
